# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") on Sheet1 is recomputed; write the new literal values for
# rows 2-35 (data rows 0-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(6,6,7,9,6,3,10,12,10,6,9,4,9,9,11,8,7,10,4,8,9,9,10,9,6,7,7,6,7,6,6,2,5,4)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
